$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'limit(sin(a2), limit(0, read(a0, a1), sub(read(a0, a1), add(write(a0, conditional(protectedDiv(0, protectedDiv(cos(a2), sin(abs(sub(cos(0), a1))))), a2), a2), a1))), protectedLog(a2))'
$ws.Cells.Item(1, 2).Value = -846
$ws.Cells.Item(2, 1).Value = 'sub(read(a0, a1), add(add(add(a2, a2), a2), add(a2, add(add(a2, write(a0, sin(0), add(sin(read(a0, limit(a2, protectedLog(a2), a2))), a2))), add(a2, sub(protectedDiv(abs(a2), a2), read(a0, read(a0, cos(conditional(a1, 0))))))))))'
$ws.Cells.Item(2, 2).Value = -435
$ws.Cells.Item(3, 1).Value = 'protectedDiv(add(read(a0, read(a0, sub(sin(add(read(a0, a1), protectedDiv(0, a1))), protectedDiv(0, a2)))), a2), conditional(sin(sin(read(a0, a1))), write(a0, limit(protectedLog(read(a0, add(add(conditional(abs(sin(a1)), a1), a2), sin(sin(a2))))), abs(limit(0, read(a0, a2), 0)), conditional(protectedDiv(0, a1), protectedDiv(a1, conditional(a2, read(a0, 0))))), a2)))'
$ws.Cells.Item(3, 2).Value = -151
$ws.Cells.Item(4, 1).Value = 'sub(sub(limit(a2, read(a0, a2), protectedLog(read(a0, abs(0)))), a2), add(write(a0, protectedLog(sub(conditional(a1, a2), sin(a2))), write(a0, protectedDiv(a1, limit(0, 0, limit(write(a0, a1, write(a0, 0, a2)), a2, a2))), a2)), limit(sin(cos(sub(protectedLog(abs(a2)), a2))), sub(0, limit(a2, a1, 0)), a2)))'
$ws.Cells.Item(4, 2).Value = -392
$ws.Cells.Item(5, 1).Value = 'sub(sub(read(a0, a1), add(write(a0, abs(a2), protectedDiv(add(a2, a2), protectedLog(write(a0, sin(cos(0)), abs(a1))))), cos(read(a0, read(a0, abs(sub(a1, conditional(limit(a2, limit(0, a1, a2), a2), a1)))))))), sub(a2, protectedLog(a2)))'
$ws.Cells.Item(5, 2).Value = -151
$ws.Cells.Item(6, 1).Value = 'sub(read(a0, a2), write(a0, sub(conditional(write(a0, a1, protectedLog(conditional(a1, protectedLog(a1)))), add(abs(0), sub(abs(write(a0, a1, a1)), add(protectedDiv(a2, a2), limit(0, 0, 0))))), conditional(0, conditional(conditional(a2, abs(a1)), a1))), protectedDiv(a2, conditional(sub(a2, sin(sin(write(a0, abs(conditional(sin(a1), 0)), a2)))), a2))))'
$ws.Cells.Item(6, 2).Value = -457
$ws.Cells.Item(7, 1).Value = 'sub(sub(0, read(a0, a1)), add(protectedDiv(conditional(a2, read(a0, read(a0, a1))), a2), add(write(a0, a1, sin(a2)), a2)))'
$ws.Cells.Item(7, 2).Value = -165
$ws.Cells.Item(8, 1).Value = 'sub(protectedDiv(sub(protectedDiv(conditional(read(a0, a2), a2), a2), sin(limit(0, 0, a1))), cos(sub(write(a0, 0, a2), cos(limit(cos(sin(a2)), 0, read(a0, 0)))))), a2)'
$ws.Cells.Item(8, 2).Value = -830
$ws.Cells.Item(9, 1).Value = 'sub(sub(read(a0, 0), protectedDiv(limit(sin(a1), protectedDiv(a2, protectedDiv(read(a0, a2), a1)), protectedLog(read(a0, a2))), a2)), write(a0, a2, a2))'
$ws.Cells.Item(9, 2).Value = -165
$ws.Cells.Item(10, 1).Value = 'protectedDiv(sub(sub(0, write(a0, a2, conditional(conditional(a2, 0), sin(read(a0, 0))))), write(a0, protectedDiv(a1, a1), conditional(conditional(a2, 0), read(a0, 0)))), a2)'
$ws.Cells.Item(10, 2).Value = -125

$ws.Range("A1:A6").RowHeight = 17.25

$ws.Range("A1:B10").Select()
